# Insert a new data row at row 126 (pushing existing rows 126:191 down to
# 127:192) and populate it with the new "Femacal de La Calera - Pepino
# ensalada" observation dated 2021-09-13 (Excel serial 44452).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("126:126").Insert()

$ws.Cells.Item(126, 1).Value  = 3
$ws.Cells.Item(126, 2).Value  = 'Femacal de La Calera'
$ws.Cells.Item(126, 3).Value  = 'Coquimbo'
$ws.Cells.Item(126, 4).Value  = 44452
$ws.Cells.Item(126, 5).Value  = 5
$ws.Cells.Item(126, 6).Value  = 100112043
$ws.Cells.Item(126, 7).Value  = 'Pepino ensalada'
$ws.Cells.Item(126, 8).Value  = 'Sin especificar'
$ws.Cells.Item(126, 9).Value  = 'Primera'
$ws.Cells.Item(126, 10).Value = 130
$ws.Cells.Item(126, 11).Value = 15500
$ws.Cells.Item(126, 12).Value = 16000
$ws.Cells.Item(126, 13).Value = 15731
$ws.Cells.Item(126, 14).Value = '$/caja 70 unidades'
$ws.Cells.Item(126, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(126, 16).Value = 225
$ws.Cells.Item(126, 17).Value = 70
$ws.Cells.Item(126, 18).Value = 'Hortaliza'
